# TC08_CDS_Filter_Platform-NovaSeq6000.xlsx — "CDS Input file updates"
#
# The ParticipantsTab Neo4j query (column B / row 2) is replaced by a revised
# version (adds OPTIONAL MATCH hops + apoc.coll.sort on the sample list), and
# every StatQuery cell (column C, rows 2-4) ends up pointing at the same
# shared-string text it already had — the reshuffle just happens because the
# old ParticipantsTab query string is dropped from the shared-string table.
# Net effect on the sheet: three cell values change and row 2 grows taller to
# fit the longer wrapped query text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$participantQuery = "MATCH (p:participant)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nOPTIONAL MATCH (samp)<--(f:file)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nWITH s, p, samp, f, g, diag`nWHERE g.platform in ['Illumina NovaSeq 6000']`nwith p`nOPTIONAL MATCH (p)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nWITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`nRETURN`ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY p.participant_id LIMIT 100"
$statQuery        = "Match (f)<--(g:genomic_info)`nWHERE g.platform in ['Illumina NovaSeq 6000']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH DISTINCT samp,s,p,f`nRETURN`n    count(distinct s) AS Studies,`n    count(distinct p) AS Participants,`n    count(distinct samp) AS Samples,`n    count(distinct f) AS ``Files``"
$sampleQuery      = "Match (f)<--(g:genomic_info)`nWHERE g.platform in ['Illumina NovaSeq 6000']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN  `n coalesce(samp.sample_id, '') as ``Sample ID``,`n coalesce(p.participant_id,'') as ``Participant ID``,`n coalesce(s.study_name, '') as ``Study Name``,`n coalesce(s.phs_accession,'') as ``Accession``,`n coalesce(samp.sample_tumor_status,'') as ``Tumor``,`ncoalesce(samp.sample_type,'') as ``Analyte Type```nORDER By samp.sample_id LIMIT 100"
$filesQuery       = "Match (f)<--(g:genomic_info)`nWHERE g.platform in ['Illumina NovaSeq 6000']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH DISTINCT p,s,samp,f`nRETURN `n    coalesce(f.file_name, '') as ``File Name``,`n    coalesce(s.study_name, '') as ``Study Name``,`n    coalesce(s.phs_accession,'') as ``Accession``,`n    coalesce(p.participant_id,'') as ``Participant ID``,`n    coalesce(samp.sample_id, '') as ``Sample ID``,`n    coalesce(f.file_type, '') as ``File Type```n   ORDER By f.file_name LIMIT 100"

# Row 2 - ParticipantsTab
$ws.Range("B2").Value = $participantQuery
$ws.Range("C2").Value = $statQuery

# Row 3 - SamplesTab
$ws.Range("B3").Value = $sampleQuery
$ws.Range("C3").Value = $statQuery

# Row 4 - FilesTab
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery

# The longer participant query now wraps across more lines, so row 2 grows.
$ws.Rows.Item(2).RowHeight = 279

# View state: scrolled down a row, new active selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("B5").Select() | Out-Null
